$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Random_Forest
$ws.Range("B3").Value = 0.96
$ws.Range("C3").Value = 0.84
$ws.Range("D3").Value = 12
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 0.71
$ws.Range("I3").Value = 0.71
$ws.Range("J3").Value = 0.71
$ws.Range("L3").Value = 0.8

# Row 6 - DNN
$ws.Range("B6").Value = 0.92
$ws.Range("C6").Value = 0.79
$ws.Range("D6").Value = 17
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 41
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.29
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0.45
$ws.Range("K6").Value = 0.09
$ws.Range("L6").Value = 0.54
